$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "에피소드#06 - Pandas 데이터 전처리, 추가, 삭제, 데이터 type 변환"
$ws.Range("E4").Value = "https://teddylee777.github.io/pandas/pandas-tutorial-06"

$ws.Range("D8").Value = "제주어 기계번역 모델과 음성합성 모델에 관한 연구를 소개합니다."
$ws.Range("E8").Value = "https://www.kakaobrain.com/blog/119"

$ws.Range("D35").Value = "[코딩 인터뷰] 화이트보드 코딩 인터뷰"
$ws.Range("E35").Value = "http://docs.likejazz.com/coding-interview-4/"

$ws.Range("D36").Value = "Towards Contrastive Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/308"

$ws.Range("D37").Value = "[Paper Review] Characteristics of DNN feature space"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1432&mod=document&pageid=1"
